# Add a new row (row 39) to the "location-1" sheet describing one more
# YouTube live-cam location (Praia de Laranjeiras, Balneario Camboriu, BR).
#
# Columns: A=Category, B=latitude/longitude, C=Location, D=CITY, E=COUNTRY,
#          F=YouTube Link, G=Status (=IsYouTubeVideoValid(F))

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new cell values first (in the same order the original author's
# data entry produced them) so the shared-string table picks up the new
# strings in the expected order.
$ws.Range("F39").Value = "StAk5P7B4Go"
$ws.Range("B39").Value = "-26.9971687101983-48.5908382934894"
$ws.Range("A39").Value = "LIVE, SEA, BEACH"
$ws.Range("C39").Value = "PRAIA DE LARANJEIRAS - BALNEÁRIO CAMBORIU - SC - RESTAURANTE CASA DO CAMARÃO - BC AO VIVO"
$ws.Range("D39").Value = "Balneário Camboriú"
$ws.Range("E39").Value = "Brazil"
$ws.Range("G39").Formula = "=IsYouTubeVideoValid(F39)"

# Copy the formatting of the row above (row 38, the previous last data row)
# down onto the new row so styles match (category/location/city/country
# cells, the quote-prefixed lat/long cell, and the status column).
$ws.Range("A38:G38").Copy()
$ws.Range("A39:G39").PasteSpecial(-4122)

# Leave the selection where the author ended up after typing the new row.
$ws.Range("G41").Select()
